$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1141
$ws.Range("F5").Value = 57
$ws.Range("F6").Value = 473
$ws.Range("F7").Value = 790
$ws.Range("F10").Value = 51
$ws.Range("F11").Value = 428
$ws.Range("F12").Value = 228
$ws.Range("F14").Value = 944
$ws.Range("F15").Value = 128
$ws.Range("F16").Value = 2055
$ws.Range("F17").Value = 538
$ws.Range("F18").Value = 8877
$ws.Range("F19").Value = 840
$ws.Range("F20").Value = 530
$ws.Range("F21").Value = 81
$ws.Range("F23").Value = 28
$ws.Range("F24").Value = 239

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 5
$ws.Range("F10").Value = 130

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5604
$ws.Range("F3").Value = 431
$ws.Range("F4").Value = 411

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5604
$ws.Range("F4").Value = 431
$ws.Range("F5").Value = 411
$ws.Range("F7").Value = 1141
$ws.Range("F8").Value = 5
$ws.Range("F10").Value = 57
$ws.Range("F11").Value = 473
$ws.Range("F12").Value = 790
$ws.Range("F16").Value = 51
$ws.Range("F17").Value = 428
$ws.Range("F18").Value = 228
$ws.Range("F22").Value = 944
$ws.Range("F24").Value = 128
$ws.Range("F25").Value = 130
$ws.Range("F27").Value = 2055
$ws.Range("F28").Value = 538
$ws.Range("F29").Value = 8877
$ws.Range("F32").Value = 840
$ws.Range("F33").Value = 530
$ws.Range("F34").Value = 81
$ws.Range("F37").Value = 28
$ws.Range("F39").Value = 239
